# Updated symbol list on Thu Dec 29 04:03:35 UTC 2022 with GitHub Actions
# Refresh the "Price" (column D) and "Hora" (column G) values for the
# crypto price table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") updates -----------------------------------------
$priceUpdates = @{
    2  = "244.91"
    3  = "23.98"
    4  = "5.203"
    6  = "6.453"
    7  = "3.210"
    8  = "0.8142"
    9  = "0.8648"
    11 = "0.06979"
    12 = "0.03190"
    13 = "0.03020"
    14 = "0.09332"
    15 = "3.822"
    16 = "0.001517"
    17 = "0.04716"
    18 = "0.0006005"
    19 = "0.006165"
    21 = "0.004106"
    22 = "0.00008693"
    23 = "3.584"
    24 = "2.144"
    27 = "0.0002326"
    40 = "0.03713"
    41 = "0.006215"
    42 = "0.1050"
    44 = "0.007547"
    45 = "0.00005226"
    47 = "0.4296"
    48 = "0.002037"
    49 = "0.00002098"
    50 = "0.0001998"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    # Keep these as plain text (matching the original inline-string cells)
    # so exact formatting like trailing zeros ("0.1050") is preserved
    # instead of being auto-coerced into a number.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# --- Column G ("Hora") updates ------------------------------------------
# Every data row (2 through 51) moves from hour "3" to hour "4".
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G$row")
    $cell.NumberFormat = "@"
    $cell.Value = "4"
}
